$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for the fixed-residue PSSM matrix (B2:K21), refreshed with supplemental-figure data
$newValues = @{
    "B2" = -18.00609592681343
    "C2" = 2.436212825300598
    "D2" = -18.00609592681343
    "E2" = -18.00609592681343
    "F2" = -18.00609592681343
    "G2" = -18.00609592681343
    "H2" = -18.00609592681343
    "I2" = -18.00609592681343
    "J2" = -18.00609592681343
    "K2" = -18.00609592681343
    "B3" = -18.00609592681343
    "C3" = -18.00609592681343
    "D3" = -18.00609592681343
    "E3" = -18.00609592681343
    "F3" = -18.00609592681343
    "G3" = -18.00609592681343
    "H3" = -18.00609592681343
    "I3" = -18.00609592681343
    "J3" = -18.00609592681343
    "K3" = -18.00609592681343
    "B4" = -18.00609592681343
    "C4" = 2.113190990022184
    "D4" = 2.826102611755927
    "E4" = -18.00609592681343
    "F4" = 2.544257257212152
    "G4" = -18.00609592681343
    "H4" = 1.938943478064196
    "I4" = -18.00609592681343
    "J4" = 2.767913303979058
    "K4" = -18.00609592681343
    "B5" = -18.00609592681343
    "C5" = 1.043648946472065
    "D5" = -18.00609592681343
    "E5" = -18.00609592681343
    "F5" = -18.00609592681343
    "G5" = 2.116303653757734
    "H5" = -18.00609592681343
    "I5" = -18.00609592681343
    "J5" = -18.00609592681343
    "K5" = -18.00609592681343
    "B6" = -18.00609592681343
    "C6" = -18.00609592681343
    "D6" = -18.00609592681343
    "E6" = -18.00609592681343
    "F6" = -18.00609592681343
    "G6" = -18.00609592681343
    "H6" = -18.00609592681343
    "I6" = -18.00609592681343
    "J6" = -18.00609592681343
    "K6" = -18.00609592681343
    "B7" = 2.964064042498262
    "C7" = -18.00609592681343
    "D7" = -18.00609592681343
    "E7" = -18.00609592681343
    "F7" = -18.00609592681343
    "G7" = -18.00609592681343
    "H7" = -18.00609592681343
    "I7" = -18.00609592681343
    "J7" = -18.00609592681343
    "K7" = -18.00609592681343
    "B8" = -18.00609592681343
    "C8" = -18.00609592681343
    "D8" = -18.00609592681343
    "E8" = 2.803107753849724
    "F8" = -18.00609592681343
    "G8" = -18.00609592681343
    "H8" = -18.00609592681343
    "I8" = -18.00609592681343
    "J8" = -18.00609592681343
    "K8" = -18.00609592681343
    "B9" = 3.608423867465102
    "C9" = -18.00609592681343
    "D9" = -18.00609592681343
    "E9" = -18.00609592681343
    "F9" = -18.00609592681343
    "G9" = -18.00609592681343
    "H9" = -18.00609592681343
    "I9" = -18.00609592681343
    "J9" = -18.00609592681343
    "K9" = -18.00609592681343
    "B10" = -18.00609592681343
    "C10" = -18.00609592681343
    "D10" = -18.00609592681343
    "E10" = -18.00609592681343
    "F10" = -18.00609592681343
    "G10" = -18.00609592681343
    "H10" = -18.00609592681343
    "I10" = -18.00609592681343
    "J10" = -18.00609592681343
    "K10" = 2.172915328965093
    "B11" = -18.00609592681343
    "C11" = -18.00609592681343
    "D11" = -18.00609592681343
    "E11" = 1.951211621109353
    "F11" = -18.00609592681343
    "G11" = 2.334728875561919
    "H11" = -18.00609592681343
    "I11" = -18.00609592681343
    "J11" = -18.00609592681343
    "K11" = 1.080188048268727
    "B12" = -18.00609592681343
    "C12" = -18.00609592681343
    "D12" = -18.00609592681343
    "E12" = -18.00609592681343
    "F12" = -18.00609592681343
    "G12" = -18.00609592681343
    "H12" = -18.00609592681343
    "I12" = -18.00609592681343
    "J12" = -18.00609592681343
    "K12" = -18.00609592681343
    "B13" = -18.00609592681343
    "C13" = -18.00609592681343
    "D13" = -18.00609592681343
    "E13" = 1.651328948309114
    "F13" = -18.00609592681343
    "G13" = -18.00609592681343
    "H13" = -18.00609592681343
    "I13" = -18.00609592681343
    "J13" = 2.061297729419057
    "K13" = 1.567686225989914
    "B14" = -18.00609592681343
    "C14" = -18.00609592681343
    "D14" = 1.64699579818443
    "E14" = -18.00609592681343
    "F14" = -18.00609592681343
    "G14" = -18.00609592681343
    "H14" = -18.00609592681343
    "I14" = -18.00609592681343
    "J14" = -18.00609592681343
    "K14" = 1.902175006002508
    "B15" = -18.00609592681343
    "C15" = -18.00609592681343
    "D15" = -0.3892527734030022
    "E15" = -18.00609592681343
    "F15" = -18.00609592681343
    "G15" = -18.00609592681343
    "H15" = -18.00609592681343
    "I15" = -18.00609592681343
    "J15" = -18.00609592681343
    "K15" = -18.00609592681343
    "B16" = -18.00609592681343
    "C16" = -18.00609592681343
    "D16" = -18.00609592681343
    "E16" = -18.00609592681343
    "F16" = -18.00609592681343
    "G16" = -18.00609592681343
    "H16" = -18.00609592681343
    "I16" = -18.00609592681343
    "J16" = 2.275065517302995
    "K16" = -18.00609592681343
    "B17" = -18.00609592681343
    "C17" = 0.483254234267301
    "D17" = -0.3201631676269855
    "E17" = -18.00609592681343
    "F17" = -18.00609592681343
    "G17" = -18.00609592681343
    "H17" = -0.1365442527538969
    "I17" = 4.321922888650185
    "J17" = 1.103278203014336
    "K17" = -18.00609592681343
    "B18" = -18.00609592681343
    "C18" = -18.00609592681343
    "D18" = -18.00609592681343
    "E18" = -18.00609592681343
    "F18" = -18.00609592681343
    "G18" = -18.00609592681343
    "H18" = -0.4140112328007287
    "I18" = -18.00609592681343
    "J18" = 1.018981437803544
    "K18" = -18.00609592681343
    "B19" = -18.00609592681343
    "C19" = -18.00609592681343
    "D19" = 1.690533767708807
    "E19" = -18.00609592681343
    "F19" = -18.00609592681343
    "G19" = -18.00609592681343
    "H19" = 1.889332913168478
    "I19" = -18.00609592681343
    "J19" = -18.00609592681343
    "K19" = -18.00609592681343
    "B20" = -18.00609592681343
    "C20" = 1.683112104036981
    "D20" = 2.317395399868282
    "E20" = -18.00609592681343
    "F20" = 3.824447519419411
    "G20" = -18.00609592681343
    "H20" = 2.406441876100954
    "I20" = -18.00609592681343
    "J20" = -18.00609592681343
    "K20" = 2.738604077121007
    "B21" = -18.00609592681343
    "C21" = 1.844287805047122
    "D21" = -18.00609592681343
    "E21" = 2.587916064665237
    "F21" = -18.00609592681343
    "G21" = 3.408651752590015
    "H21" = 2.459216378619995
    "I21" = -18.00609592681343
    "J21" = -18.00609592681343
    "K21" = -18.00609592681343
}

foreach ($cellRef in $newValues.Keys) {
    $ws.Range($cellRef).Value = $newValues[$cellRef]
}
